# cementFactory.xlsx - recycle flows usable at factory level
# Rebuilds the "Connections" sheet with new Origin_Unit / Destination_Unit /
# Origin_Chain / Destination_Chain / Recycle_Replacing / Purge_Fraction /
# Max_Replace_Fraction columns, updates the cementFactoryConnections defined
# name to start at column B, and makes "Connections" the active sheet/tab.

$wb = $excel.ActiveWorkbook

$wsConnections = $wb.Worksheets.Item("Connections")
$wsPower       = $wb.Worksheets.Item("Power Chain")
$wsCO2Capture  = $wb.Worksheets.Item("CO2 Capture")

# --- Rebuild the Connections table (new column layout) -------------------
# New column order:
#  A Product | B Origin_Chain | C Origin_Unit | D Product_IO_of_Origin |
#  E Destination_Chain | F Product_IO_of_Destination | G Recycle_Replacing |
#  H Destination_Unit | I Purge_Fraction | J Max_Replace_Fraction

$headers = @("Product","Origin_Chain","Origin_Unit","Product_IO_of_Origin","Destination_Chain","Product_IO_of_Destination","Recycle_Replacing","Destination_Unit","Purge_Fraction","Max_Replace_Fraction")
for ($col = 1; $col -le $headers.Length; $col++) {
    $wsConnections.Cells.Item(1, $col).Value = $headers[$col - 1]
}

$rows = @(
    @("CO2","cement","kiln","outflow","CO2capture","inflow",$null,$null,$null,$null),
    @("electricity","cement","all","inflow","power","outflow",$null,$null,$null,$null),
    @("electricity","CO2capture","all","inflow","power","outflow",$null,$null,$null,$null),
    @("waste heat","cement","kiln","outflow","power","inflow","fuel","PowerStation",0.2,$null),
    @("compressedCO2","CO2capture","CO2Compression","outflow","cement","inflow","gypsum","blender",$null,0.1)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $rowData = $rows[$i]
    for ($col = 1; $col -le $rowData.Length; $col++) {
        $val = $rowData[$col - 1]
        if ($null -ne $val) {
            $wsConnections.Cells.Item($r, $col).Value = $val
        }
    }
}

# --- Defined name: cementFactoryConnections now starts at column B -------
$wb.Names.Item("cementFactoryConnections").RefersTo = "=Connections!`$B`$1:`$F`$4"

# --- Selections on the other chain sheets ---------------------------------
$wsPower.Range("B2").Select()
$wsCO2Capture.Range("B3").Select()

# --- Make Connections the active sheet/tab, selection at A5 --------------
$wsConnections.Range("A5").Select()
$wsConnections.Activate()
